$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating the "2021-Q4" sheet (same
#    layout/styling) and placing it right before the "总计" sheet.
# ---------------------------------------------------------------------------
$srcQuarter = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$srcQuarter.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Fill in the new quarter's fund data (row 2). The fund code/name are
# unchanged from the template ("2021-Q4") so they are left as-is. Columns
# D/E/F/G must stay text (leading apostrophe forces text storage, matching
# the source sheets) while A2 and H2 remain numeric.
$newSheet.Cells.Item(2,4).Value = "'0.48"
$newSheet.Cells.Item(2,5).Value = "'89.68"
$newSheet.Cells.Item(2,6).Value = "'4.59"
$newSheet.Cells.Item(2,7).Value = "'0.0220"
$newSheet.Cells.Item(2,8).Value = 1

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row for "2022-Q1" at the
#    top of the data (row 2) and push the existing rows down.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# The inserted row picks up stray formatting on B:D - clear it so it matches
# the unstyled data cells used by the other rows.
$ws.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index cells in column A (copy format
# from A3, which still holds the original row-2 formatting).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "2022-Q1"
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = 0.02

# Renumber the index column for the rows that shifted down.
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(7,1).Value = 5
